$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of Argent (Solar) price data appended as row 95.
# Source values are stored as plain text in the workbook (matching the
# existing rows), so force the Text number format before writing the
# values to prevent Excel from auto-converting them to dates/numbers.
$rng = $ws.Range("A95:J95")
$rng.NumberFormat = "@"

$ws.Range("A95").Value = "2025-06-04"
$ws.Range("B95").Value = "35.5"
$ws.Range("C95").Value = "35.01"
$ws.Range("D95").Value = "0.94"
$ws.Range("E95").Value = "0.253"
$ws.Range("F95").Value = "0.09"
$ws.Range("G95").Value = "5,501"
$ws.Range("H95").Value = "8,236"
$ws.Range("I95").Value = "8,286"
$ws.Range("J95").Value = "7.2166"
